# The "optimization_parameters" sheet had a stray leftover row labelled
# "Sheet" (with values 3 / 4) sitting between the "Strain" row and the
# "simulation_timepoints" row. Remove it - this shifts every row below it
# up by one and drops the now-unused "Sheet" shared string / number-format
# style automatically.
$wb = $excel.ActiveWorkbook

$paramsSheet = $wb.Worksheets.Item("optimization_parameters")
$paramsSheet.Rows(16).Delete()

# The author finished auditing on the last sheet (optimization_diagnostics),
# so that tab ends up the active one when the file is saved.
$diagSheet = $wb.Worksheets.Item("optimization_diagnostics")
$diagSheet.Activate()
